$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell, forcing text storage (with original
# no-style formatting) even when the text looks like a number, matching
# the source data which stores prices/volumes as literal text strings.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    if ($text -match "^\s*[+-]?(\d+\.?\d*|\.\d+)\s*$") {
        $cell.Formula = "'" + $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

# Row 37/38: Hedera and Algorand swap ranking positions, with updated price/volume data
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D37" "0.2099"
Set-TextValue "E37" "  +0.68%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D38" "0.06045"
Set-TextValue "E38" "  +1.24%  "

# Remaining numeric/percentage updates
Set-TextValue "D2" '27.582.60'
Set-TextValue "E2" '  -0.62%  '
Set-TextValue "D3" '1.751.28'
Set-TextValue "E3" '  +0.52%  '
Set-TextValue "D4" '1.004'
Set-TextValue "D5" '323.95'
Set-TextValue "E5" '  +1.32%  '
Set-TextValue "E6" '  -0.23%  '
Set-TextValue "D7" '0.4599'
Set-TextValue "E7" '  +9.47%  '
Set-TextValue "E8" '  -1.24%  '
Set-TextValue "D9" '0.07491'
Set-TextValue "E9" '  +1.66%  '
Set-TextValue "D10" '42.13'
Set-TextValue "E10" '  -1.57%  '
Set-TextValue "E11" '  +1.55%  '
Set-TextValue "E12" '  -0.28%  '
Set-TextValue "D13" '20.77'
Set-TextValue "E13" '  +1.26%  '
Set-TextValue "D14" '6.009'
Set-TextValue "E14" '  -0.22%  '
Set-TextValue "D15" '7.084'
Set-TextValue "E15" '  -1.93%  '
Set-TextValue "D16" '1.754.18'
Set-TextValue "E16" '  -1.56%  '
Set-TextValue "D17" '92.48'
Set-TextValue "E17" '  +2.11%  '
Set-TextValue "E18" '  +1.55%  '
Set-TextValue "D19" '0.06431'
Set-TextValue "E19" '  +1.45%  '
Set-TextValue "D20" '1.002'
Set-TextValue "E20" '  -0.10%  '
Set-TextValue "D21" '16.79'
Set-TextValue "E21" '  -0.46%  '
Set-TextValue "D22" '5.810'
Set-TextValue "E22" '  -1.53%  '
Set-TextValue "D23" '27.648.77'
Set-TextValue "E23" '  -0.46%  '
Set-TextValue "E24" '  +0.57%  '
Set-TextValue "D25" '2.107'
Set-TextValue "E25" '  +1.78%  '
Set-TextValue "D26" '164.20'
Set-TextValue "E26" '  +4.95%  '
Set-TextValue "E27" '  +1.77%  '
Set-TextValue "D28" '1.955.06'
Set-TextValue "E28" '  -1.19%  '
Set-TextValue "D29" '2.073'
Set-TextValue "E29" '  -2.05%  '
Set-TextValue "D30" '126.46'
Set-TextValue "E30" '  +2.57%  '
Set-TextValue "E31" '  -5.66%  '
Set-TextValue "D32" '0.09193'
Set-TextValue "E32" '  +4.43%  '
Set-TextValue "D33" '3.664'
Set-TextValue "E33" '  +0.63%  '
Set-TextValue "D34" '5.532'
Set-TextValue "E34" '  +0.18%  '
Set-TextValue "D35" '11.88'
Set-TextValue "E35" '  -2.43%  '
Set-TextValue "D36" '0.02295'
Set-TextValue "E36" '  +1.60%  '
Set-TextValue "D39" '4.978'
Set-TextValue "E39" '  +1.34%  '
Set-TextValue "D40" '0.6330'
Set-TextValue "E40" '  +1.29%  '
Set-TextValue "D41" '1.210'
Set-TextValue "E41" '  +3.34%  '
Set-TextValue "D42" '1.379'
Set-TextValue "E42" '  -0.95%  '
Set-TextValue "D43" '7.783'
Set-TextValue "E43" '  +0.10%  '
Set-TextValue "D44" '13.33'
Set-TextValue "E44" '  +0.37%  '
Set-TextValue "D45" '0.5915'
Set-TextValue "E45" '  +1.55%  '
Set-TextValue "D46" '3.713'
Set-TextValue "E46" '  +1.20%  '
Set-TextValue "D47" '122.96'
Set-TextValue "E47" '  +1.03%  '
Set-TextValue "E48" '  -0.66%  '
Set-TextValue "D49" '1.141'
Set-TextValue "E49" '  -2.07%  '
Set-TextValue "D50" '0.06859'
Set-TextValue "E50" '  +1.05%  '
Set-TextValue "D51" '72.14'
Set-TextValue "E51" '  -1.25%  '
